# Weekly refresh of fruit/vegetable prices (Haba, Macroferia Regional de Talca).
# The underlying records got reshuffled and three new weekly observations
# were appended. Rather than replaying the shuffle, just drive every row
# (existing 2-12 plus the new 13-15) to its final target content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that stay constant for every data row already present in the sheet.
$constA = 5
$constB = "Macroferia Regional de Talca"
$constC = "Maule"
$constE = 7
$constF = 100112026
$constG = "Haba"
$constH = "Sin especificar"
$constI = "Primera"
$constN = "`$/saco 25 kilos"
$constQ = 25
$constR = "Hortaliza"

# Final per-row values for D (fecha), J (volumen), K (precio minimo),
# L (precio maximo), M (precio promedio ponderado), O (origen), P (precio $/Kg).
$rows = @(
    @{ Row = 2;  D = 44383; J = 120; K = 12000; L = 12000; M = 12000; O = "Provincia del Elquí"; P = 480 },
    @{ Row = 3;  D = 44445; J = 200; K = 10000; L = 10000; M = 10000; O = "Provincia del Elquí"; P = 400 },
    @{ Row = 4;  D = 44166; J = 200; K = 8000;  L = 8000;  M = 8000;  O = "Región del Maule";    P = 320 },
    @{ Row = 5;  D = 44162; J = 200; K = 9000;  L = 9000;  M = 9000;  O = "Región del Maule";    P = 360 },
    @{ Row = 6;  D = 44165; J = 150; K = 7000;  L = 7000;  M = 7000;  O = "Región del Maule";    P = 280 },
    @{ Row = 7;  D = 44169; J = 200; K = 9000;  L = 9000;  M = 9000;  O = "Región del Maule";    P = 360 },
    @{ Row = 8;  D = 44159; J = 300; K = 7000;  L = 7000;  M = 7000;  O = "Región del Maule";    P = 280 },
    @{ Row = 9;  D = 44396; J = 200; K = 14000; L = 14000; M = 14000; O = "Provincia del Elquí"; P = 560 },
    @{ Row = 10; D = 44168; J = 200; K = 9000;  L = 9000;  M = 9000;  O = "Región del Maule";    P = 360 },
    @{ Row = 11; D = 44398; J = 200; K = 15000; L = 15000; M = 15000; O = "Provincia del Elquí"; P = 600 },
    @{ Row = 12; D = 44441; J = 200; K = 10000; L = 10000; M = 10000; O = "Provincia del Elquí"; P = 400 },
    @{ Row = 13; D = 44446; J = 200; K = 10000; L = 10000; M = 10000; O = "Provincia del Elquí"; P = 400 },
    @{ Row = 14; D = 44161; J = 200; K = 9000;  L = 9000;  M = 9000;  O = "Región del Maule";    P = 360 },
    @{ Row = 15; D = 44167; J = 200; K = 8000;  L = 8000;  M = 8000;  O = "Región del Maule";    P = 320 }
)

# Existing dimension only covered rows 1-12; rows 13-15 are brand new and need
# every column (not just the ones that vary) populated.
$dateNumberFormat = $ws.Range("D2").NumberFormat

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Cells.Item($r, 1).Value2 = $constA
    $ws.Cells.Item($r, 2).Value2 = $constB
    $ws.Cells.Item($r, 3).Value2 = $constC
    $ws.Cells.Item($r, 4).Value2 = $item.D
    $ws.Cells.Item($r, 4).NumberFormat = $dateNumberFormat
    $ws.Cells.Item($r, 5).Value2 = $constE
    $ws.Cells.Item($r, 6).Value2 = $constF
    $ws.Cells.Item($r, 7).Value2 = $constG
    $ws.Cells.Item($r, 8).Value2 = $constH
    $ws.Cells.Item($r, 9).Value2 = $constI
    $ws.Cells.Item($r, 10).Value2 = $item.J
    $ws.Cells.Item($r, 11).Value2 = $item.K
    $ws.Cells.Item($r, 12).Value2 = $item.L
    $ws.Cells.Item($r, 13).Value2 = $item.M
    $ws.Cells.Item($r, 14).Value2 = $constN
    $ws.Cells.Item($r, 15).Value2 = $item.O
    $ws.Cells.Item($r, 16).Value2 = $item.P
    $ws.Cells.Item($r, 17).Value2 = $constQ
    $ws.Cells.Item($r, 18).Value2 = $constR
}
